$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 33   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/2/2026  Through  2/8/2026"

# --- Table updates (rows 14-33) ---
# Stable donor cells used purely to carry over cell *style* (via PasteSpecial -> Formats only)
# after a Value assignment flips a cell between "text placeholder" (style 13) and "numeric" (style 14/15).
# None of these donors are themselves edited anywhere in this script.
$styleDonorText = $ws.Range("C14")   # style 13 (text placeholder, e.g. "0" / "***.*")
$styleDonorCount = $ws.Range("J14")  # style 14 (plain integer count)
$styleDonorPct = $ws.Range("K14")    # style 15 (percent-change number)

# Row 14
$ws.Range("G14").Value = "'0"
$styleDonorText.Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = "'***.*"
$styleDonorText.Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null

# Row 15
$ws.Range("G15").Value = "'0"
$styleDonorText.Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = "'***.*"
$styleDonorText.Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("N15").Value = -75

# Row 16
$ws.Range("C16").Value = 1
$styleDonorCount.Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = -53.846153846153
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = -81.818181818181
$ws.Range("N16").Value = -95.081967213114

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 26
$ws.Range("K17").Value = -3.846153846153
$ws.Range("L17").Value = 92.307692307692
$ws.Range("M17").Value = 127.272727272727
$ws.Range("N17").Value = -47.916666666666

# Row 18
$ws.Range("C18").Value = 1
$styleDonorCount.Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 4
$ws.Range("K18").Value = -42.857142857142
$ws.Range("L18").Value = -63.636363636363
$ws.Range("M18").Value = -66.666666666666
$ws.Range("N18").Value = -94.666666666666

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = -18.181818181818
$ws.Range("I19").Value = 26
$ws.Range("J19").Value = 29
$ws.Range("K19").Value = -10.344827586206
$ws.Range("L19").Value = 30
$ws.Range("M19").Value = -18.75
$ws.Range("N19").Value = -49.019607843137

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "'0"
$styleDonorText.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "'***.*"
$styleDonorText.Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 7
$ws.Range("K20").Value = -22.222222222222
$ws.Range("L20").Value = -65
$ws.Range("M20").Value = -12.5
$ws.Range("N20").Value = -91.463414634146

# Row 21
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = -13.333333333333
$ws.Range("F21").Value = 47
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = -18.965517241379
$ws.Range("I21").Value = 69
$ws.Range("J21").Value = 88
$ws.Range("K21").Value = -21.590909090909
$ws.Range("L21").Value = -5.479452054794
$ws.Range("M21").Value = -28.125
$ws.Range("N21").Value = -82.03125

# Row 22
$ws.Range("C22").Value = 2
$styleDonorCount.Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = 2
$styleDonorCount.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = 0
$styleDonorPct.Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$styleDonorCount.Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null
$ws.Range("H22").Value = 0
$styleDonorPct.Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 2
$styleDonorCount.Copy() | Out-Null
$ws.Range("J22").PasteSpecial(-4122) | Out-Null
$ws.Range("K22").Value = 50
$styleDonorPct.Copy() | Out-Null
$ws.Range("K22").PasteSpecial(-4122) | Out-Null
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -50

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = -33.333333333333
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -18.181818181818
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = -35.294117647058
$ws.Range("L23").Value = 175
$ws.Range("M23").Value = -8.333333333333

# Row 24
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 33
$ws.Range("H24").Value = -15.384615384615
$ws.Range("I24").Value = 50
$ws.Range("J24").Value = 55
$ws.Range("K24").Value = -9.090909090909
$ws.Range("L24").Value = -3.846153846153
$ws.Range("M24").Value = -47.368421052631

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = "'0"
$styleDonorText.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = "'***.*"
$styleDonorText.Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 44.444444444444
$ws.Range("I25").Value = 16
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = -5.882352941176

# Row 26
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 22
$ws.Range("H26").Value = 10
$ws.Range("I26").Value = 28
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = 12
$ws.Range("L26").Value = -20
$ws.Range("M26").Value = 21.739130434782

# Row 27
$ws.Range("G27").Value = "'0"
$styleDonorText.Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").Value = "'***.*"
$styleDonorText.Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null

# Row 28
$ws.Range("C28").Value = "'0"
$styleDonorText.Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = -80

# Row 33
$ws.Range("C33").Value = "'0"
$styleDonorText.Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0